# Weekly fruit/vegetable price update for Brócoli @ Vega Central Mapocho de Santiago.
# - Rows 339-344 get new/corrected price readings (some dates, qualities and
#   price columns change) and rows 345-348 get a date correction
#   (44400 -> 44399).
# - Four brand-new rows (349-352) are appended for 2021-07-23 (serial 44400).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 339 ----
$ws.Range("D339").Value = 44448
$ws.Range("J339").Value = 4300
$ws.Range("K339").Value = 600
$ws.Range("L339").Value = 650
$ws.Range("M339").Value = 625
$ws.Range("P339").Value = 625

# ---- Row 340 ----
$ws.Range("D340").Value = 44448
$ws.Range("J340").Value = 5200
$ws.Range("K340").Value = 600
$ws.Range("L340").Value = 650
$ws.Range("M340").Value = 625
$ws.Range("O340").Value = "Región de O'Higgins"
$ws.Range("P340").Value = 625

# ---- Row 341 ----
$ws.Range("D341").Value = 44448
$ws.Range("I341").Value = "Segunda"
$ws.Range("J341").Value = 1600
$ws.Range("K341").Value = 450
$ws.Range("L341").Value = 500
$ws.Range("M341").Value = 475
$ws.Range("P341").Value = 475

# ---- Row 342 ----
$ws.Range("D342").Value = 44448
$ws.Range("I342").Value = "Segunda"
$ws.Range("J342").Value = 2100
$ws.Range("K342").Value = 450
$ws.Range("L342").Value = 500
$ws.Range("M342").Value = 475
$ws.Range("P342").Value = 475

# ---- Row 343 ----
$ws.Range("D343").Value = 44167
$ws.Range("I343").Value = "Primera"
$ws.Range("J343").Value = 4500
$ws.Range("K343").Value = 500
$ws.Range("L343").Value = 600
$ws.Range("M343").Value = 544
$ws.Range("P343").Value = 544

# ---- Row 344 ----
$ws.Range("D344").Value = 44238
$ws.Range("I344").Value = "Primera"
$ws.Range("J344").Value = 3400
$ws.Range("K344").Value = 700
$ws.Range("L344").Value = 800
$ws.Range("M344").Value = 750
$ws.Range("O344").Value = "Región Metropolitana"
$ws.Range("P344").Value = 750

# ---- Rows 345-348: date correction only (44400 -> 44399) ----
$ws.Range("D345").Value = 44399
$ws.Range("D346").Value = 44399
$ws.Range("D347").Value = 44399
$ws.Range("D348").Value = 44399

# ---- New rows 349-352 ----
$dateFmt = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A349").Value = 9
$ws.Range("B349").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C349").Value = "Metropolitana"
$ws.Range("D349").Value = 44400
$ws.Range("D349").NumberFormat = $dateFmt
$ws.Range("E349").Value = 13
$ws.Range("F349").Value = 100112023
$ws.Range("G349").Value = "Brócoli"
$ws.Range("H349").Value = "Sin especificar"
$ws.Range("I349").Value = "Primera"
$ws.Range("J349").Value = 5200
$ws.Range("K349").Value = 500
$ws.Range("L349").Value = 600
$ws.Range("M349").Value = 550
$ws.Range("N349").Value = "`$/unidad"
$ws.Range("O349").Value = "Región Metropolitana"
$ws.Range("P349").Value = 550
$ws.Range("Q349").Value = 1
$ws.Range("R349").Value = "Hortaliza"

$ws.Range("A350").Value = 9
$ws.Range("B350").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C350").Value = "Metropolitana"
$ws.Range("D350").Value = 44400
$ws.Range("D350").NumberFormat = $dateFmt
$ws.Range("E350").Value = 13
$ws.Range("F350").Value = 100112023
$ws.Range("G350").Value = "Brócoli"
$ws.Range("H350").Value = "Sin especificar"
$ws.Range("I350").Value = "Primera"
$ws.Range("J350").Value = 4300
$ws.Range("K350").Value = 500
$ws.Range("L350").Value = 600
$ws.Range("M350").Value = 550
$ws.Range("N350").Value = "`$/unidad"
$ws.Range("O350").Value = "Región de O'Higgins"
$ws.Range("P350").Value = 550
$ws.Range("Q350").Value = 1
$ws.Range("R350").Value = "Hortaliza"

$ws.Range("A351").Value = 9
$ws.Range("B351").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C351").Value = "Metropolitana"
$ws.Range("D351").Value = 44400
$ws.Range("D351").NumberFormat = $dateFmt
$ws.Range("E351").Value = 13
$ws.Range("F351").Value = 100112023
$ws.Range("G351").Value = "Brócoli"
$ws.Range("H351").Value = "Sin especificar"
$ws.Range("I351").Value = "Segunda"
$ws.Range("J351").Value = 2500
$ws.Range("K351").Value = 400
$ws.Range("L351").Value = 400
$ws.Range("M351").Value = 400
$ws.Range("N351").Value = "`$/unidad"
$ws.Range("O351").Value = "Región Metropolitana"
$ws.Range("P351").Value = 400
$ws.Range("Q351").Value = 1
$ws.Range("R351").Value = "Hortaliza"

$ws.Range("A352").Value = 9
$ws.Range("B352").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C352").Value = "Metropolitana"
$ws.Range("D352").Value = 44400
$ws.Range("D352").NumberFormat = $dateFmt
$ws.Range("E352").Value = 13
$ws.Range("F352").Value = 100112023
$ws.Range("G352").Value = "Brócoli"
$ws.Range("H352").Value = "Sin especificar"
$ws.Range("I352").Value = "Segunda"
$ws.Range("J352").Value = 1600
$ws.Range("K352").Value = 400
$ws.Range("L352").Value = 400
$ws.Range("M352").Value = 400
$ws.Range("N352").Value = "`$/unidad"
$ws.Range("O352").Value = "Región de O'Higgins"
$ws.Range("P352").Value = 400
$ws.Range("Q352").Value = 1
$ws.Range("R352").Value = "Hortaliza"

Write-Host "Done updating rows 339-352"
